$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Schedule": insert one new pump-session row (old 12h block split into
# a 4h block + a new 8h block), shifting the final row down.
# ---------------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Schedule")

# Insert a fresh row above the current row 2; old row2->row3, old row3->row4.
$wsSchedule.Rows.Item(2).Insert()
# The inserted row inherits the header's bold style by default - reset it.
$wsSchedule.Rows.Item(2).Style = "Normal"
# Columns A & B hold datetimes in this sheet; copy the datetime format down
# from the row below so the new row renders the same way as its neighbours.
$wsSchedule.Range("A2:B2").NumberFormat = $wsSchedule.Range("A3:B3").NumberFormat

# ---------------------------------------------------------------------------
# Sheet "Detailed": insert one new half-hour sample at the top of the series,
# shifting every following row down by one.
# ---------------------------------------------------------------------------
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsDetailed.Rows.Item(2).Insert()
$wsDetailed.Rows.Item(2).Style = "Normal"
# Column A holds a datetime, column D a date - copy down the formats used by
# the rest of the table.
$wsDetailed.Range("A2").NumberFormat = $wsDetailed.Range("A3").NumberFormat
$wsDetailed.Range("D2").NumberFormat = $wsDetailed.Range("D3").NumberFormat

# ---------------------------------------------------------------------------
# Rewrite the data for both tables (run 137 output).
# ---------------------------------------------------------------------------
$sheet1Data = @(
  @(46043, 46043.16666666666, 4, 15.12, 523.8023115, 34.64301001984127),
  @(46043.33333333334, 46043.66666666666, 8, 30.24, -238.73240625, -7.89459015376984),
  @(46043.83333333334, 46044, 4, 15.12, 489.610875, 32.38167162698413)
)

$sheet2Data = @(
  @(46043, 64.8901, "historical", 46043, "ON"),
  @(46043.02083333334, 64.89, "historical", 46043, "ON"),
  @(46043.04166666666, 64.89, "historical", 46043, "ON"),
  @(46043.0625, 67.82252, "historical", 46043, "ON"),
  @(46043.08333333334, 65, "historical", 46043, "ON"),
  @(46043.10416666666, 65, "forecast", 46043, "ON"),
  @(46043.125, 66.80052000000001, "forecast", 46043, "ON"),
  @(46043.14583333334, 77.94, "forecast", 46043, "ON"),
  @(46043.16666666666, 78, "forecast", 46043, "OFF"),
  @(46043.1875, 78.14865, "forecast", 46043, "OFF"),
  @(46043.20833333334, 87.63102000000001, "forecast", 46043, "OFF"),
  @(46043.22916666666, 101.25, "forecast", 46043, "OFF"),
  @(46043.25, 101.25, "forecast", 46043, "OFF"),
  @(46043.27083333334, 77.94, "forecast", 46043, "OFF"),
  @(46043.29166666666, 56.82404, "forecast", 46043, "OFF"),
  @(46043.3125, 0.7, "forecast", 46043, "OFF"),
  @(46043.33333333334, 0, "forecast", 46043, "ON"),
  @(46043.35416666666, -5.62975, "forecast", 46043, "ON"),
  @(46043.375, -6.47731, "forecast", 46043, "ON"),
  @(46043.39583333334, -6.77623, "forecast", 46043, "ON"),
  @(46043.41666666666, -8.94598, "forecast", 46043, "ON"),
  @(46043.4375, -9.5, "forecast", 46043, "ON"),
  @(46043.45833333334, -14, "forecast", 46043, "ON"),
  @(46043.47916666666, -15.67332, "forecast", 46043, "ON"),
  @(46043.5, -16.16716, "forecast", 46043, "ON"),
  @(46043.52083333334, -22.35389, "forecast", 46043, "ON"),
  @(46043.54166666666, -22.06114, "forecast", 46043, "ON"),
  @(46043.5625, -23.5, "forecast", 46043, "ON"),
  @(46043.58333333334, -23.5, "forecast", 46043, "ON"),
  @(46043.60416666666, -23.5, "forecast", 46043, "ON"),
  @(46043.625, -23.5, "forecast", 46043, "ON"),
  @(46043.64583333334, -23.26897, "forecast", 46043, "ON"),
  @(46043.66666666666, -6.8, "forecast", 46043, "OFF"),
  @(46043.6875, -5.51, "forecast", 46043, "OFF"),
  @(46043.70833333334, 36.06, "forecast", 46043, "OFF"),
  @(46043.72916666666, 48.31538, "forecast", 46043, "OFF"),
  @(46043.75, 55.63438, "forecast", 46043, "OFF"),
  @(46043.77083333334, 64.8901, "forecast", 46043, "OFF"),
  @(46043.79166666666, 77.94, "forecast", 46043, "OFF"),
  @(46043.8125, 79.95, "forecast", 46043, "OFF"),
  @(46043.83333333334, 73.29000000000001, "forecast", 46043, "ON"),
  @(46043.85416666666, 65, "forecast", 46043, "ON"),
  @(46043.875, 57.35973, "forecast", 46043, "ON"),
  @(46043.89583333334, 60.01917, "forecast", 46043, "ON"),
  @(46043.91666666666, 57.09, "forecast", 46043, "ON"),
  @(46043.9375, 61.93302, "forecast", 46043, "ON"),
  @(46043.95833333334, 63.52943, "forecast", 46043, "ON"),
  @(46043.97916666666, 63.94365, "forecast", 46043, "ON")
)

$r = 2
foreach ($row in $sheet1Data) {
  $wsSchedule.Cells.Item($r, 1).Value = $row[0]
  $wsSchedule.Cells.Item($r, 2).Value = $row[1]
  $wsSchedule.Cells.Item($r, 3).Value = $row[2]
  $wsSchedule.Cells.Item($r, 4).Value = $row[3]
  $wsSchedule.Cells.Item($r, 5).Value = $row[4]
  $wsSchedule.Cells.Item($r, 6).Value = $row[5]
  $r = $r + 1
}

$r = 2
foreach ($row in $sheet2Data) {
  $wsDetailed.Cells.Item($r, 1).Value = $row[0]
  $wsDetailed.Cells.Item($r, 2).Value = $row[1]
  $wsDetailed.Cells.Item($r, 3).Value = $row[2]
  $wsDetailed.Cells.Item($r, 4).Value = $row[3]
  $wsDetailed.Cells.Item($r, 5).Value = $row[4]
  $r = $r + 1
}
